# send intermediate results to cloud server
# Adds new "server_t4_latency" / "server_gpu_latency" / "server_cpu_latency"
# benchmark columns, relabels the existing latency columns, and marks a
# staging cell (J1) with a highlight fill.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row relabel / additions -----------------------------------
# E was "server_latency" -> now "pi_latency"
# F was "pi_latency"     -> now "server_gpu_latency"
# G (new)                -> "server_t4_latency"
# H (new)                -> "server_cpu_latency"
$ws.Range("E1").Value = "pi_latency"
$ws.Range("F1").Value = "server_gpu_latency"
$ws.Range("G1").Value = "server_t4_latency"
$ws.Range("H1").Value = "server_cpu_latency"

# Highlighted placeholder cell for the new cloud-server upload marker
$ws.Range("J1").Interior.Color = 65535

# --- New "server_t4_latency" (G) / "server_cpu_latency" (H) data ------
$ws.Range("G2").Value = 0.0001
$ws.Range("H2").Value = 0.0001

$ws.Range("G3").Value = 0.0015999999999999973
$ws.Range("H3").Value = 0.0044866197183098516

$ws.Range("G4").Value = 0.001606153846153843
$ws.Range("H4").Value = 0.00515985915492957

$ws.Range("G5").Value = 0.0028515384615384635
$ws.Range("H5").Value = 0.021180281690140825

$ws.Range("G6").Value = 0.004033846153846158
$ws.Range("H6").Value = 0.03910845070422537

$ws.Range("G7").Value = 0.005421538461538462
$ws.Range("H7").Value = 0.06307535211267605

$ws.Range("G8").Value = 0.007065384615384625
$ws.Range("H8").Value = 0.07645704225352108

$ws.Range("G9").Value = 0.007123076923076932
$ws.Range("H9").Value = 0.06658169014084508

$ws.Range("G10").Value = 0.007178461538461546
$ws.Range("H10").Value = 0.057695070422535265

# --- Move the active selection to match the saved cursor position -----
$ws.Range("G19").Select()
